$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped from 45190 (2023-09-21)
# to 45192 (2023-09-23) for every data row (rows 2 through 511).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 511 }

$ws.Range("C2:C$lastRow").Value = 45192
